{"js": "// Fix the typo in the document title: \"eCoacing Log System\" -> \"eCoaching Log System\"\n// (a missing \"h\" is inserted so \"eCoacing\" reads \"eCoaching\").\nconst body = context.document.body;\n\n// Search for the misspelled word; match case so we don't disturb any already-correct\n// \"eCoaching\" occurrences elsewhere in the document.\nconst results = body.search(\"eCoacing\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the text 'eCoacing' to fix.\");\n}\n\n// Replace the misspelled word with the corrected spelling. insertText(..., \"Replace\")\n// swaps just the matched text in place, keeping the run's existing character\n// formatting (bold, size, etc.) untouched.\nresults.items.forEach((range) => {\n  range.insertText(\"eCoaching\", Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Fix the typo in the document title: \"eCoacing Log System\" -> \"eCoaching Log System\"\n# (a missing \"h\" is inserted so \"eCoacing\" reads \"eCoaching\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"eCoacing\"\n$find.Replacement.Text = \"eCoaching\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n\n# wdReplaceAll = 2 -- replace every match (there is exactly one: the title).\n$find.Execute([ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n"}
